$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-2024")

$ws.Range("A27").Value = 45370
$ws.Range("B27").Value = "MPAL"
$ws.Range("C27").Value = "TP"
$ws.Range("D27").Value = "x"
$ws.Range("G27").Value = "ScrabbleScore"
$ws.Range("I27").Value = "présentation de la stratégie TDD à employer, puis travail en autonomie"

$ws.Range("A28").Value = 45370
$ws.Range("B28").Value = "MPAL"
$ws.Range("C28").Value = "TP"
$ws.Range("F28").Value = "x"
$ws.Range("G28").Value = "ScrabbleScore"
$ws.Range("I28").Value = "présentation de la stratégie TDD à employer, puis travail en autonomie"
